$d = $word.ActiveDocument

$p18 = $d.Paragraphs.Item(18)
$p18.Range.Delete()

$delRange = $d.Range($d.Paragraphs.Item(22).Range.Start, $d.Paragraphs.Item(29).Range.End)
$delRange.Delete()

Write-Output ("Content.End=" + $d.Content.End)
$p21 = $d.Paragraphs.Item(21)
Write-Output ("p21 End=" + $p21.Range.End)

$bmPoint1 = $d.Range(2845, 2845)
$d.Bookmarks.Add("TESTA", $bmPoint1)
Write-Output "added TESTA at 2845"

$bmPoint2 = $d.Range(2846, 2846)
$d.Bookmarks.Add("TESTB", $bmPoint2)
Write-Output "added TESTB at 2846"

$bmPoint3 = $d.Range(2800, 2800)
$d.Bookmarks.Add("TESTC", $bmPoint3)
Write-Output "added TESTC at 2800"
